# ACTL3142 DATA.xlsx - "Add files via upload" commit replay
#
# The commit adds a second worksheet ("Sheet1") after the existing
# "ACTL3142 DATA" sheet. The new sheet re-lays the monthly claims table
# (previously at B3:G18 on "ACTL3142 DATA") with the months running down
# column A and the years across columns B:F, adds a new
# "avg. claims/quarter" summary row, and becomes the active/selected tab.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Add the new worksheet right after "ACTL3142 DATA" ------------------
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$new.Name = "Sheet1"

# --- Header row: year labels (first literal, rest formulas) -------------
$new.Range("B1").Value = 2018
$new.Range("C1").Formula = "=1+B1"
$new.Range("D1").Formula = "=1+C1"
$new.Range("E1").Formula = "=1+D1"
$new.Range("F1").Formula = "=1+E1"
$new.Range("A1:F1").Font.Bold = $true
$new.Range("A1:F1").Font.Size = 12

# --- Month labels (column A) --------------------------------------------
$new.Range("A2").Value = "Jan"
$new.Range("A3").Value = "Feb"
$new.Range("A4").Value = "Mar"
$new.Range("A5").Value = "Apr"
$new.Range("A6").Value = "May"
$new.Range("A7").Value = "Jun"
$new.Range("A8").Value = "Jul"
$new.Range("A9").Value = "Aug"
$new.Range("A10").Value = "Sep"
$new.Range("A11").Value = "Oct"
$new.Range("A12").Value = "Nov"
$new.Range("A13").Value = "Dec"

# --- Monthly claim counts, 2018-2022 (2022 only has Jan-Jun data) -------
$new.Range("B2").Value = 559
$new.Range("C2").Value = 845
$new.Range("D2").Value = 824
$new.Range("E2").Value = 764
$new.Range("F2").Value = 527

$new.Range("B3").Value = 704
$new.Range("C3").Value = 901
$new.Range("D3").Value = 876
$new.Range("E3").Value = 1009
$new.Range("F3").Value = 810

$new.Range("B4").Value = 910
$new.Range("C4").Value = 1019
$new.Range("D4").Value = 993
$new.Range("E4").Value = 1039
$new.Range("F4").Value = 969

$new.Range("B5").Value = 928
$new.Range("C5").Value = 957
$new.Range("D5").Value = 614
$new.Range("E5").Value = 986
$new.Range("F5").Value = 789

$new.Range("B6").Value = 1092
$new.Range("C6").Value = 1157
$new.Range("D6").Value = 667
$new.Range("E6").Value = 1103
$new.Range("F6").Value = 911

$new.Range("B7").Value = 1029
$new.Range("C7").Value = 971
$new.Range("D7").Value = 857
$new.Range("E7").Value = 1039
$new.Range("F7").Value = 494

$new.Range("B8").Value = 1031
$new.Range("C8").Value = 1144
$new.Range("D8").Value = 950
$new.Range("E8").Value = 790

$new.Range("B9").Value = 1065
$new.Range("C9").Value = 1106
$new.Range("D9").Value = 911
$new.Range("E9").Value = 555

$new.Range("B10").Value = 885
$new.Range("C10").Value = 946
$new.Range("D10").Value = 932
$new.Range("E10").Value = 501

$new.Range("B11").Value = 965
$new.Range("C11").Value = 1095
$new.Range("D11").Value = 985
$new.Range("E11").Value = 543

$new.Range("B12").Value = 1045
$new.Range("C12").Value = 969
$new.Range("D12").Value = 989
$new.Range("E12").Value = 826

$new.Range("B13").Value = 798
$new.Range("C13").Value = 899
$new.Range("D13").Value = 875
$new.Range("E13").Value = 774

# --- "total" row ----------------------------------------------------------
$new.Range("A14").Value = "total"
$new.Range("B14").Formula = "=SUM(B2:B13)"
$new.Range("C14").Formula = "=SUM(C2:C13)"
$new.Range("D14").Formula = "=SUM(D2:D13)"
$new.Range("E14").Formula = "=SUM(E2:E13)"
$new.Range("F14").Formula = "=SUM(F2:F7)"

# --- "avg. claims/quarter" row (new row introduced by this commit) ------
$new.Range("A15").Value = "avg. claims/quarter"
$new.Range("B15").Formula = "=B14/4"
$new.Range("C15").Formula = "=C14/4"
$new.Range("D15").Formula = "=D14/4"
$new.Range("E15").Formula = "=E14/4"
$new.Range("F15").Formula = "=F14/2"

# --- "avg. claims/month" row ----------------------------------------------
$new.Range("A16").Value = "avg. claims/month"
$new.Range("B16").Formula = "=B14/12"
$new.Range("C16").Formula = "=C14/12"
$new.Range("D16").Formula = "=D14/12"
$new.Range("E16").Formula = "=E14/12"
$new.Range("F16").Formula = "=F14/6"

# --- Column A width (label column) ---------------------------------------
$new.Columns.Item(1).ColumnWidth = 15

# --- View state: new sheet is the active/selected tab --------------------
$new.Range("D10").Select()

# --- View state on the original data sheet: selection now A1:H18 --------
$dataSheet.Range("A1:H18").Select()

# Re-activate the new sheet last so it is the workbook's active tab
# (matches activeTab="1" / tabSelected on the new sheet in the target file).
$new.Activate()
